# GripperBOM.xlsx update: rename/extend "Electronic parts" section with new
# connector rows, add a "Sample links" column of reference URLs, widen col B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Other mechanical parts" table: add a sample link for the cylindrical rollers row ---
$ws.Range("F19").Value = "Sample links"
$ws.Range("F19").Font.Bold = $true

$ws.Range("F21").Value = "https://www.ebmia.pl/waleczki/37040-walek-waleczki-5x8-skf.html"

# --- "Electronic parts" table: sample links for existing rows ---
$ws.Range("F28").Value = "https://botland.store/standard-servos/3576-servo-powerhd-lf-20mg-standard-6939670200387.html"
$ws.Range("F29").Value = "https://www.sparkfun.com/products/15100"

# Row 31 used to just say "HZ03" with no amount; rename + give it an amount + link
$ws.Range("A31").Value = "HZ03 Connector with pins"
$ws.Range("D31").Value = 1
$ws.Range("F31").Value = "https://en.maritex.com.pl/connectors/crimp_terminal_connectors/crimp_terminal_2_54mm_pitch/female_crimp_terminal_housings_for_cable_with_snap-lock_2_54_mm_pitch/hz03.html"

# New row 32: HZ02 connector
$ws.Range("A32:B32").Merge()
$ws.Range("A32").Value = "HZ02 Connector with pins"
$ws.Range("A32").HorizontalAlignment = -4108
$ws.Range("B32").HorizontalAlignment = -4108
$ws.Range("D32").Value = 1
$ws.Range("F32").Value = "https://en.maritex.com.pl/connectors/crimp_terminal_connectors/crimp_terminal_2_54mm_pitch/female_crimp_terminal_housings_for_cable_with_snap-lock_2_54_mm_pitch/hz02.html"

# New row 33: USB-C male plug
$ws.Range("A33:B33").Merge()
$ws.Range("A33").Value = "USB-C male plug(to solder wires)"
$ws.Range("A33").HorizontalAlignment = -4108
$ws.Range("B33").HorizontalAlignment = -4108
$ws.Range("D33").Value = 1
$ws.Range("F33").Value = "https://www.amazon.com/Solder-Cable-Socket-Attached-Board/dp/B07P1BDNQV"

# New row 34: Wires (no amount)
$ws.Range("A34:B34").Merge()
$ws.Range("A34").Value = "Wires"
$ws.Range("A34").HorizontalAlignment = -4108
$ws.Range("B34").HorizontalAlignment = -4108

# Widen column B to fit the longer item names now in the sheet
$ws.Columns(2).ColumnWidth = 25.8333333333333

# Restore the selection to roughly where the author left off editing
$ws.Range("O41").Select()
